$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "renters" (sheet1.xml): add new "Owner Credit" column (Q)
# ---------------------------------------------------------------------------
$renters = $wb.Worksheets.Item("renters")

$renters.Range("Q1").Value = "Owner Credit"
$renters.Range("Q2:Q9").Value = 700

# ---------------------------------------------------------------------------
# Sheet "owners" (sheet2.xml): flip "Owner Insurance" off, update some
# damage states, and add the new "Owner Credit" column (P)
# ---------------------------------------------------------------------------
$owners = $wb.Worksheets.Item("owners")

$owners.Range("C2:C9").Value = 0

$owners.Range("L2").Value = "Extensive"
$owners.Range("L5").Value = "Moderate"

$owners.Range("P1").Value = "Owner Credit"
$owners.Range("P2").Value = 500
$owners.Range("P3:P9").Value = 700

# ---------------------------------------------------------------------------
# Restore the selections shown for each sheet (also refreshes dimensions)
# ---------------------------------------------------------------------------
$renters.Range("M22").Select()
$owners.Range("L6").Select()
